$d = $word.ActiveDocument

# Locate the last paragraph in the document body (the one ending
# "...ja yrittäisin tehdä siitä, niin hyvän, kuin vain voi."), and
# append a brand-new paragraph right after it, before the sectPr.
$lastIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($lastIndex)
$lastPara.Range.InsertParagraphAfter()

$newIndex = $d.Paragraphs.Count
$newPara = $d.Paragraphs.Item($newIndex)
$newRange = $newPara.Range

$newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:left="2603" w:hanging="2603"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Kehitys</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:tab/><w:t xml:space="preserve">Jatkokehitystä on aloitettu, muttei ole vielä saatu täysin valmiiksi. </w:t></w:r></w:p>'

$newRange.InsertXML($newParaXml) | Out-Null
